# Fruta / hortaliza, semanal
# A new weekly reading is prepended as row 4 (Fecha 2021-11-19 / serial 44503),
# pushing the previously-existing rows (old rows 4..22) down by one (new rows 5..23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4; Excel shifts rows 4..22 down to 5..23
# and the used-range dimension grows from R22 to R23 automatically.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with this week's market reading.
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44503
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15429
$ws.Range("N4").Value = "`$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1187
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
